$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (GitHub Actions scheduled refresh).
# D-column "Price" values are forced to Text (leading apostrophe) so strings
# like "1.50" / "0.0000255" keep their exact digits instead of Excel
# auto-coercing them into numbers (which would drop trailing zeros / use
# scientific notation). B/C/E columns are already non-numeric-looking text.

$ws.Range("D2").Value = "'95.337.30"
$ws.Range("E2").Value = "  -1.19%  "
$ws.Range("D3").Value = "'3.587.38"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "'236.99"
$ws.Range("E5").Value = "  -1.99%  "
$ws.Range("D6").Value = "'657.48"
$ws.Range("E6").Value = "  +1.77%  "
$ws.Range("D7").Value = "'1.50"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "'0.402"
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("E9").Value = "  +0.17%  "
$ws.Range("D10").Value = "'1.02"
$ws.Range("E10").Value = "  +0.96%  "
$ws.Range("D11").Value = "'3.586.82"
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("D13").Value = "'42.72"
$ws.Range("E13").Value = "  -3.18%  "
$ws.Range("D14").Value = "'6.49"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "'4.259.44"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").Value = "'95.511.14"
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "'0.0000255"
$ws.Range("E17").Value = "  -1.00%  "
$ws.Range("D18").Value = "'3.580.59"
$ws.Range("E18").Value = "  -2.49%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'12.88"
$ws.Range("E19").Value = "  -4.71%  "
$ws.Range("B20").Value = "Polkadot"
$ws.Range("C20").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D20").Value = "'7.77"
$ws.Range("E20").Value = "  -8.16%  "
$ws.Range("D21").Value = "'17.98"
$ws.Range("E21").Value = "  -2.88%  "
$ws.Range("D22").Value = "'3.49"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'0.490"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "'511.31"
$ws.Range("E24").Value = "  -1.81%  "
$ws.Range("D25").Value = "'7.09"
$ws.Range("E25").Value = "  +4.30%  "
$ws.Range("D26").Value = "'0.0000198"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").Value = "'95.49"
$ws.Range("E27").Value = "  -2.53%  "
$ws.Range("D28").Value = "'12.91"
$ws.Range("E28").Value = "  +2.28%  "
$ws.Range("D29").Value = "'3.781.90"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("D30").Value = "'3.06"
$ws.Range("E30").Value = "  -3.76%  "
$ws.Range("E31").Value = "  +1.83%  "
$ws.Range("D32").Value = "'11.63"
$ws.Range("E32").Value = "  -0.86%  "
$ws.Range("E33").Value = "  +0.30%  "
$ws.Range("D34").Value = "'1.00"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "'0.178"
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").Value = "'32.27"
$ws.Range("E36").Value = "  +0.89%  "
$ws.Range("D37").Value = "'1.71"
$ws.Range("E37").Value = "  +14.31%  "
$ws.Range("D38").Value = "'8.73"
$ws.Range("E38").Value = "  +10.29%  "
$ws.Range("D39").Value = "'0.564"
$ws.Range("E39").Value = "  -2.73%  "
$ws.Range("D40").Value = "'599.92"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("E41").Value = "  -0.48%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("B43").Value = "ImmutableX"
$ws.Range("C43").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D43").Value = "'1.90"
$ws.Range("E43").Value = "  +8.95%  "
$ws.Range("B44").Value = "ARBITRUM"
$ws.Range("C44").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D44").Value = "'0.917"
$ws.Range("E44").Value = "  -3.03%  "
$ws.Range("D45").Value = "'35.78"
$ws.Range("E45").Value = "  +18.18%  "
$ws.Range("D46").Value = "'5.78"
$ws.Range("E46").Value = "  -0.38%  "
$ws.Range("E47").Value = "  +2.98%  "
$ws.Range("D48").Value = "'23.42"
$ws.Range("E48").Value = "  -1.52%  "
$ws.Range("D49").Value = "'0.0418"
$ws.Range("E49").Value = "  -3.02%  "
$ws.Range("E50").Value = "  +0.02%  "
$ws.Range("D51").Value = "'8.23"
$ws.Range("E51").Value = "  -0.41%  "
